# Fix: Auto-focus description field when editing a single-item sales
# transaction, not just multi-item.
#
# Data-model side effect of the fix: the placeholder "Reference"/"Actions"
# cells that used to be written for the previously-active row (14) are no
# longer there once a new row becomes the active one; a new single-item
# transaction row (15) is appended to the Sales sheet instead, carrying
# those placeholder cells now.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales")

# Row 14 is no longer the "active" row, so drop its empty placeholder cells.
$ws.Range("D14").ClearContents()
$ws.Range("H14").ClearContents()

# Append the new single-item transaction as row 15.
$ws.Range("A15").Value = "تجربة عميل4"

# Plain text date - force text entry with a leading apostrophe (otherwise it
# gets auto-converted to a date serial number), then reset the cell's style
# back to Normal so no extra formatting ends up attached to the cell.
$bCell = $ws.Range("B15")
$bCell.Value = "'2025-07-18"
$bCell.Style = "Normal"

$ws.Range("C15").Value = "#1: بروش | Qty: 2 | Price: 20 | Total: 40 | VAT: 6"

# Empty text placeholder cells (Reference / Actions columns), same pattern
# used previously by row 14 while it was the active row.
$dCell = $ws.Range("D15")
$dCell.Value = "'"
$dCell.Style = "Normal"

$ws.Range("E15").Value = 40
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 46

$hCell = $ws.Range("H15")
$hCell.Value = "'"
$hCell.Style = "Normal"

$ws.Range("I15").Value = $false
